$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the report title (shared string) to the corrected Kyrgyz wording.
$ws.Range("A1").Value = "8.10.2.2 Камсыздандыруу компаниялардын финансылык көрсөткүчтөрү"

# 2. Add the new "2023" column (Q) of data, mirroring column P's formatting
#    by copying the source cell (which also copies its style) and then
#    overwriting the value with the new figure for 2023.
$ws.Range("P3").Copy($ws.Range("Q3"))
$ws.Range("Q3").Value = 2023

$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("Q4").Value = 16

$ws.Range("P5").Copy($ws.Range("Q5"))
$ws.Range("Q5").Value = 3031.4

# 3. Restore the default selection at A1 (the saved view previously pointed
#    at E9, outside of the data range).
$ws.Range("A1").Select()
